$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.101.67'
$ws.Range("E2").Value = '  -3.93%  '

$ws.Range("D3").Value = '1.860.23'
$ws.Range("E3").Value = '  -4.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.08'
$ws.Range("E5").Value = '  -4.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4664'
$ws.Range("E7").Value = '  -3.30%  '

$ws.Range("E8").Value = '  -4.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06546'
$ws.Range("E9").Value = '  -4.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.63'
$ws.Range("E10").Value = '  -2.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07818'
$ws.Range("E11").Value = '  -0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.43'
$ws.Range("E12").Value = '  -8.17%  '

$ws.Range("D13").Value = '1.858.34'
$ws.Range("E13").Value = '  -4.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.123'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6651'
$ws.Range("E15").Value = '  -4.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '280.48'
$ws.Range("E16").Value = '  -6.22%  '

$ws.Range("D17").Value = '30.133.07'
$ws.Range("E17").Value = '  -3.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  +0.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.503'
$ws.Range("E19").Value = '  -2.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.57'
$ws.Range("E20").Value = '  -3.82%  '

$ws.Range("D21").Value = '2.104.56'
$ws.Range("E21").Value = '  -4.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007222'
$ws.Range("E22").Value = '  -5.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.123'
$ws.Range("E24").Value = '  -5.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.310'
$ws.Range("E25").Value = '  -3.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.58'
$ws.Range("E26").Value = '  -1.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.85'
$ws.Range("E27").Value = '  -5.73%  '

$ws.Range("E28").Value = '  -11.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.334'
$ws.Range("E29").Value = '  -4.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09543'
$ws.Range("E30").Value = '  -6.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.401'
$ws.Range("E31").Value = '  -5.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.465'
$ws.Range("E32").Value = '  -4.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.095'
$ws.Range("E33").Value = '  -6.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04645'
$ws.Range("E34").Value = '  -4.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7013'
$ws.Range("E35").Value = '  -6.30%  '

$ws.Range("E36").Value = '  -4.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.695'
$ws.Range("E37").Value = '  -1.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01852'
$ws.Range("E38").Value = '  -5.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.275'
$ws.Range("E39").Value = '  -5.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.507'
$ws.Range("E40").Value = '  -5.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.86'
$ws.Range("E41").Value = '  -5.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8533'
$ws.Range("E42").Value = '  -3.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.914'
$ws.Range("E43").Value = '  -6.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  +0.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.78'
$ws.Range("E45").Value = '  -2.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4145'
$ws.Range("E46").Value = '  -5.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '994.16'
$ws.Range("E47").Value = '  -2.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.179'
$ws.Range("E48").Value = '  -5.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.290'
$ws.Range("E49").Value = '  +0.91%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.03'
$ws.Range("E50").Value = '  -3.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1137'
$ws.Range("E51").Value = '  -6.96%  '
